$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores plain-text values in the source workbook.
# Some replacement values (e.g. "238.09", "0.07867") would otherwise be
# auto-detected by Excel as numbers, so force a Text number format on those
# specific cells first to make sure they remain text, matching the source data.
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.324.16"
$ws.Range("E2").Value = "  +0.37%  "

$ws.Range("D3").Value = "1.860.79"
$ws.Range("E3").Value = "  +0.07%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D6").Value = "238.09"
$ws.Range("E6").Value = "  +0.32%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").Value = "0.07867"
$ws.Range("E8").Value = "  +1.68%  "

$ws.Range("D9").Value = "0.3053"
$ws.Range("E9").Value = "  +0.25%  "

$ws.Range("D10").Value = "24.93"
$ws.Range("E10").Value = "  +7.20%  "

$ws.Range("D11").Value = "2.027.83"
$ws.Range("E11").Value = "  +8.74%  "

$ws.Range("E12").Value = "  +0.07%  "

$ws.Range("D13").Value = "5.227"
$ws.Range("E13").Value = "  +1.33%  "

$ws.Range("D14").Value = "0.7192"
$ws.Range("E14").Value = "  +0.10%  "

$ws.Range("D15").Value = "89.66"
$ws.Range("E15").Value = "  +0.62%  "

$ws.Range("D16").Value = "29.330.19"
$ws.Range("E16").Value = "  +0.36%  "

$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").Value = "5.848"
$ws.Range("E17").Value = "  +1.41%  "

$ws.Range("B18").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C18").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D18").Value = "2.292.70"
$ws.Range("E18").Value = "  +8.91%  "

$ws.Range("D19").Value = "0.000007807"
$ws.Range("E19").Value = "  +0.97%  "

$ws.Range("E20").Value = "  -0.61%  "

$ws.Range("D21").Value = "238.49"
$ws.Range("E21").Value = "  +0.71%  "

$ws.Range("D22").Value = "0.9998"
$ws.Range("E22").Value = "  +0.00%  "

$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.07%  "

$ws.Range("D24").Value = "7.590"
$ws.Range("E24").Value = "  +1.95%  "

$ws.Range("D25").Value = "162.86"
$ws.Range("E25").Value = "  +0.36%  "

$ws.Range("D26").Value = "8.919"
$ws.Range("E26").Value = "  -0.97%  "

$ws.Range("D27").Value = "0.1434"
$ws.Range("E27").Value = "  -3.34%  "

$ws.Range("D28").Value = "18.10"
$ws.Range("E28").Value = "  +0.54%  "

$ws.Range("D29").Value = "1.904"
$ws.Range("E29").Value = "  -6.61%  "

$ws.Range("D30").Value = "1.376"
$ws.Range("E30").Value = "  -3.94%  "

$ws.Range("D31").Value = "1.477"
$ws.Range("E31").Value = "  -0.22%  "

$ws.Range("D32").Value = "4.327"
$ws.Range("E32").Value = "  -2.43%  "

$ws.Range("D33").Value = "4.058"
$ws.Range("E33").Value = "  +0.65%  "

$ws.Range("D34").Value = "0.05200"

$ws.Range("D35").Value = "1.180"
$ws.Range("E35").Value = "  +0.95%  "

$ws.Range("D36").Value = "0.7162"
$ws.Range("E36").Value = "  +1.25%  "

$ws.Range("D37").Value = "1.009"
$ws.Range("E37").Value = "  +0.85%  "

$ws.Range("E38").Value = "  +0.32%  "

$ws.Range("D39").Value = "0.01857"
$ws.Range("E39").Value = "  +0.64%  "

$ws.Range("D40").Value = "2.691"
$ws.Range("E40").Value = "  -1.25%  "

$ws.Range("D41").Value = "1.176.82"
$ws.Range("E41").Value = "  +3.09%  "

$ws.Range("D42").Value = "0.9221"
$ws.Range("E42").Value = "  -1.36%  "

$ws.Range("D43").Value = "6.028"
$ws.Range("E43").Value = "  +2.15%  "

$ws.Range("D44").Value = "71.89"
$ws.Range("E44").Value = "  +1.55%  "

$ws.Range("D45").Value = "0.4277"
$ws.Range("E45").Value = "  +0.10%  "

$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "0.9998"
$ws.Range("E46").Value = "  -0.02%  "

$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "2.170.48"
$ws.Range("E47").Value = "  +8.30%  "

$ws.Range("D48").Value = "102.24"
$ws.Range("E48").Value = "  -1.07%  "

$ws.Range("D49").Value = "0.5330"
$ws.Range("E49").Value = "  -2.20%  "

$ws.Range("D50").Value = "1.767"
$ws.Range("E50").Value = "  -1.46%  "

$ws.Range("D51").Value = "9.196"
$ws.Range("E51").Value = "  +0.26%  "
